$wb = $excel.ActiveWorkbook

# --- Update descriptive text on the "OPIS" sheet ---
$opis = $wb.Worksheets.Item("OPIS")
$opis.Range("B5").Value = "Średnia cena za 1 m2 lokali mieszkalnych sprzedanych w ramach transakcji rynkowych"
$opis.Range("B6").Value = "Uzyskana jako iloraz wartości i powierzchni użytkowej lokali mieszkalnych zaliczonych do danego grupowania. Cena obejmuje wszystkie składniki nieruchomości objętej transakcją. Ceny lokali mieszkalnych sprzedanych na rynku pierwotnym zawierają podatek VAT."

# --- Update column F values (price per m2) on the "DANE" sheet ---
$dane = $wb.Worksheets.Item("DANE")

$newValues = @{
    2 = 4093
    3 = 4886
    4 = 4111
    5 = 4678
    6 = 8407
    7 = 9587
    8 = 3001
    9 = 3672
    10 = 5700
    11 = 6246
    12 = 3684
    13 = 3844
    14 = 6580
    15 = 7065
    16 = 4651
    17 = 5169
    18 = 4265
    19 = 5466
    20 = 4153
    21 = 4792
    22 = 7178
    23 = 8030
    24 = 5184
    25 = 5690
    26 = 4745
    27 = 5403
    28 = 4924
    29 = 6189
    30 = 5750
    31 = 6489
    32 = 4449
    33 = 5326
    34 = 5532
    35 = 6053
    36 = 9339
    37 = 10923
    38 = 4818
    39 = 5240
    40 = 4394
    41 = 5249
    42 = 4721
    43 = 5470
    44 = 4805
    45 = 5730
    46 = 5031
    47 = 5800
    48 = 5860
    49 = 5824
    50 = 11775
    51 = 13235
    52 = 5744
    53 = 6381
    54 = 4065
    55 = 4498
    56 = 3487
    57 = 4019
    58 = 6269
    59 = 7326
    60 = 3787
    61 = 4554
    62 = 6593
    63 = 7602
    64 = 5238
    65 = 6117
    66 = 5371
    67 = 6455
    68 = 9808
    69 = 11025
    70 = 9021
    71 = 9754
    72 = 4368
    73 = 5210
    74 = 13929
    75 = 15107
    76 = 5276
    77 = 6017
    78 = 2981
    79 = 3425
    80 = 4162
    81 = 4474
    82 = 4800
    83 = 5351
    84 = 4180
    85 = 4828
    86 = 4668
    87 = 5305
    88 = 3309
    89 = 3642
    90 = 4539
    91 = 4927
    92 = 5863
    93 = 6488
    94 = 3934
    95 = 4770
    96 = 3476
    97 = 3819
    98 = 3613
    99 = 4090
    100 = 3726
    101 = 4298
    102 = 4396
    103 = 4640
    104 = 4130
    105 = 4586
    106 = 3524
    107 = 4230
    108 = 5510
    109 = 6009
    110 = 3675
    111 = 4044
    112 = 4204
    113 = 5049
    114 = 5684
    115 = 6592
    116 = 4994
    117 = 5561
    118 = 6206
    119 = 6851
    120 = 4178
    121 = 5344
    122 = 3918
    123 = 4763
    124 = 4042
    125 = 4640
    126 = 7371
    127 = 8598
    130 = 6362
    131 = 7314
    132 = 8751
    133 = 11131
}

foreach ($row in $newValues.Keys) {
    $dane.Cells.Item($row, 6).Value = $newValues[$row]
}
